# The deck's active design ("Integral") theme (ppt/theme/theme2.xml) is
# swapped with the unused default "Office Theme" (ppt/theme/theme1.xml):
# the colour scheme that used to render the slides/slideMaster becomes the
# plain Office colours, while the Office colour set moves into the spare
# theme part. Font scheme and format scheme are identical between the two
# theme parts already, so only the 12 theme colours need to change.

function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# -- this is the scheme that becomes the presentation's live colour scheme.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $tcs.Colors($i).RGB = HexToRgbInt $officeColors[$i - 1]
}
